$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value  = 201.7009458245258
$ws.Range("C5").Value  = 60.68899876601913
$ws.Range("C9").Value  = 17.40995165602148
$ws.Range("C11").Value = -94.44648335722724
$ws.Range("C12").Value = 39.96873956131847
$ws.Range("C13").Value = -61.01363796799191
$ws.Range("C15").Value = 41.58579920735638
$ws.Range("C16").Value = -10.48954692552599
$ws.Range("C20").Value = -733.8373775254573
$ws.Range("C24").Value = -414.9014438867886
$ws.Range("C26").Value = -444.4839643649241
$ws.Range("C27").Value = -514.3560478476167
$ws.Range("C28").Value = -253.3854986589261
$ws.Range("C29").Value = -477.0838833015364
$ws.Range("C32").Value = -369.5371911757881
